# This workbook is a weekly price log: every week a brand new record is
# inserted at the top of the data block (row 11, right after the most
# recent existing record in row 10), pushing all the older records down
# by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11 (shifts existing rows 11..84 down to 12..85).
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C11").Value = "Metropolitana"
$ws.Range("D11").Value = 45163
$ws.Range("E11").Value = 13
$ws.Range("F11").Value = 100112010
$ws.Range("G11").Value = "Achicoria"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 97
$ws.Range("K11").Value = 7000
$ws.Range("L11").Value = 7000
$ws.Range("M11").Value = 7000
$ws.Range("N11").Value = "`$/caja 16 unidades"
$ws.Range("O11").Value = "Provincia de Quillota"
$ws.Range("P11").Value = 438
$ws.Range("Q11").Value = 16
$ws.Range("R11").Value = "Hortaliza"
